$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update metrics for year 2025 row (row 6)
$ws.Range("C6").Value = 410
$ws.Range("E6").Value = 103
$ws.Range("G6").Value = 25.12195121951219
$ws.Range("H6").Value = 74.8780487804878
